$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.095.82'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '1.788.22'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.57'
$ws.Range("E5").Value = '  +0.95%  '

$ws.Range("E6").Value = '  -0.60%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.15'
$ws.Range("E8").Value = '  -1.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.295'
$ws.Range("E9").Value = '  +2.95%  '

$ws.Range("E10").Value = '  -3.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'
$ws.Range("E11").Value = '  +0.71%  '

$ws.Range("D12").Value = '2.045.68'
$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.33'
$ws.Range("E13").Value = '  +4.08%  '

$ws.Range("D14").Value = '1.781.65'
$ws.Range("E14").Value = '  -0.77%  '

$ws.Range("D15").Value = '34.064.72'
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.621'
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.17'
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.78'
$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("E20").Value = '  -1.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.94'
$ws.Range("E21").Value = '  +1.84%  '

$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("E24").Value = '  -3.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.84'
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("E26").Value = '  +1.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.28'
$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("E29").Value = '  +0.29%  '

$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0518'
$ws.Range("E31").Value = '  +0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.61'
$ws.Range("E33").Value = '  +2.75%  '

$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("D35").Value = '1.403.66'
$ws.Range("E35").Value = '  +0.64%  '

$ws.Range("E36").Value = '  +0.48%  '

$ws.Range("E37").Value = '  -0.70%  '

$ws.Range("E38").Value = '  +1.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.34'
$ws.Range("E39").Value = '  +5.38%  '

$ws.Range("E40").Value = '  +1.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '80.04'
$ws.Range("E41").Value = '  +1.32%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.917'
$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.71'
$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.27'
$ws.Range("E44").Value = '  +8.65%  '

$ws.Range("E45").Value = '  -6.91%  '

$ws.Range("E46").Value = '  +2.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.03'
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("E48").Value = '  +0.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.70'
$ws.Range("E49").Value = '  -1.31%  '

$ws.Range("D50").Value = '1.946.97'
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("E51").Value = '  +0.26%  '
